$d = $word.ActiveDocument

function Set-ParaXml($paragraph, [string]$innerXml) {
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($wrapper)
}

# --- Paragraph 4: "Add the script..." -> "Create the ISummonable Spells" (+ moved bookmark) ---
$p4 = $d.Paragraphs(4)
$inner4 = '<w:r><w:t xml:space="preserve">Create the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ISummonable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Spells</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParaXml $p4 $inner4

# --- Insert two new paragraphs after paragraph 6 ("Create a mask...") ---
$p6 = $d.Paragraphs(6)
$p6.Range.InsertParagraphAfter()
$pEmpty = $d.Paragraphs(7)
Set-ParaXml $pEmpty ''

$pEmpty.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(8)
$inner8 = '<w:r><w:t>1/2/2016 The seal recognizes when the player selected a valid pattern and after all that validation destroys the spell seal.</w:t></w:r>'
Set-ParaXml $pNew $inner8

# --- Paragraph 9 (was 7): "1/1/2017 ..." -- add proofErr wraps, drop bookmark (moved to p4) ---
$p9 = $d.Paragraphs(9)
$inner9 = '<w:r><w:t xml:space="preserve">1/1/2017 </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">The spell Seal registers input only from the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>summoner' + [char]0x2019 + 's</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> hand and if it does not have a spell in it. </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Fixed a mistake in the teleportation indicator. The </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>lazerPOinter</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> script was looking for extra conditions to indicate that the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>lazer</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> was leaving the platform. In fact, the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>lazer</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> was not leaving. I just took out the last conditions and it worked fine. Now you just have to check if there was a previous contact and whether the hit transform is null. If that happens, means that the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>lazer</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> left the platform.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>This is still not quite what I want because I need to check whether the platform is empty, but that is good enough for now.</w:t></w:r>'
Set-ParaXml $p9 $inner9

# --- Paragraph 13 (was 11): "12/23/2016 ... SpellInfos." ---
$p13 = $d.Paragraphs(13)
$inner13 = '<w:r><w:t xml:space="preserve">12/23/2016 Fixed the fireball feel. Added the shield. Refactored the spells. Now we do not need </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SpellInfos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>'
Set-ParaXml $p13 $inner13

# --- Paragraph 15 (was 13): "12/10/2016 ... Expandible ... spellcaster ..." ---
$p15 = $d.Paragraphs(15)
$inner15 = '<w:r><w:t xml:space="preserve">12/10/2016 Modified the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Expandible</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> class. Now spells will be thrown by the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>spellcaster</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> instead of moving by themselves.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> Fixed the state machine because you could not have spells in both hands. Added the fireball spell.</w:t></w:r>'
Set-ParaXml $p15 $inner15

# --- Paragraph 16 (was 14): "12/3/2016 ... Honovi ..." ---
$p16 = $d.Paragraphs(16)
$inner16 = '<w:r><w:t>12/3/2016 Changed the teleportation to the right touchpad press</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">. Added the state machine in </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Honovi</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>. Now the state handles the input from the player.</w:t></w:r>'
Set-ParaXml $p16 $inner16

# --- Paragraph 19 (was 17): "11/14/2016 ... Honovi ..." ---
$p19 = $d.Paragraphs(19)
$inner19 = '<w:r><w:t xml:space="preserve">11/14/2016 Created a new class for </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Honovi</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> that inherits from the</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> abstract class Character. The g</w:t></w:r>' + `
  '<w:r><w:t>ame is still working.</w:t></w:r>'
Set-ParaXml $p19 $inner19

# --- Paragraph 21 (was 19): "11/11/2016 ... Honovi ... IOManager ...") ---
$p21 = $d.Paragraphs(21)
$inner21 = '<w:r><w:t>11/11/2016</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Added the JSON File for </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Honovi</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> and the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>IOManager</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> that will send it to the player class, which will pass it to the character. </w:t></w:r>'
Set-ParaXml $p21 $inner21

# --- Paragraph 23 (was 21): "11/5/2016 ... Honovi info ..." ---
$p23 = $d.Paragraphs(23)
$inner23 = '<w:r><w:t>11/5/2016</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">I started the implementation of the Character class. It currently has dummy data with </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Honovi</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> info. You have to call the class from Player. And implement the reactions. </w:t></w:r>'
Set-ParaXml $p23 $inner23

# --- Paragraph 25 (was 23): "8/29/2016 ... Raycast ..." ---
$p25 = $d.Paragraphs(25)
$inner25 = '<w:r><w:t xml:space="preserve">8/29/2016 </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">The Prototype of the game is teleporting using </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Raycast</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>. Removed the ability to jump to avoid motion sickness.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> Added the ability to throw balls.</w:t></w:r>'
Set-ParaXml $p25 $inner25

Write-Output "done"
